# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on the BSM and GSM leve-profit tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# BSM sheet: newly populated H:N columns for rows that previously had no
# market-board data (most rows get H,I,J,K,L with M/N left blank unless the
# diff also supplied a profit figure).
# ---------------------------------------------------------------------------
$wsBSM = $wb.Worksheets.Item("BSM")

$bsmRows = @{
    117 = @{ H = 0;         I = 0;        J = 0;         K = 0;        L = 0 }
    118 = @{ H = 0;         I = 0;        J = 0;         K = 0;        L = 0 }
    119 = @{ H = 35000;     I = 0;        J = 35000;     K = 0;        L = 35000;     N = -44676 }
    120 = @{ H = 0;         I = 0;        J = 0;         K = 0;        L = 0 }
    122 = @{ H = 49775;     I = 0;        J = 49775;     K = 0;        L = 49775;     N = -59575 }
    123 = @{ H = 53703.53;  I = 0;        J = 53703.53;  K = 0;        L = 53703.53;  N = -63503.53 }
    124 = @{ H = 30000;     I = 0;        J = 30000;     K = 0;        L = 30000;     N = -39820 }
    125 = @{ H = 54980;     I = 0;        J = 54980;     K = 0;        L = 54980;     N = -64820 }
    126 = @{ H = 54980;     I = 0;        J = 54980;     K = 0;        L = 54980;     N = -64860 }
    127 = @{ H = 43966.8;   I = 0;        J = 43966.8;   K = 0;        L = 43966.8;   N = -53886.8 }
    128 = @{ H = 1800;      I = 1800;     J = 0;         K = 5400;     L = 0;         M = -2910 }
    129 = @{ H = 49989.5;   I = 0;        J = 49989.5;   K = 0;        L = 49989.5;   N = -59989.5 }
    130 = @{ H = 415495.12; I = 0;        J = 415495.12; K = 0;        L = 415495.12; N = -425535.12 }
    131 = @{ H = 38780;     I = 0;        J = 38780;     K = 0;        L = 38780;     N = -48860 }
    132 = @{ H = 44314.668; I = 0;        J = 44314.668; K = 0;        L = 44314.668; N = -54434.668 }
    133 = @{ H = 40580;     I = 0;        J = 40580;     K = 0;        L = 40580;     N = -50700 }
    134 = @{ H = 22323560;  I = 29413428; J = 5105312;   K = 88240284; L = 15315936;  M = -88237749; N = -15321006 }
    135 = @{ H = 39900;     I = 0;        J = 39900;     K = 0;        L = 39900;     N = -50040 }
    137 = @{ H = 44780;     I = 0;        J = 44780;     K = 0;        L = 44780;     N = -54980 }
    138 = @{ H = 43266.668; I = 0;        J = 43266.668; K = 0;        L = 43266.668; N = -53546.668 }
    139 = @{ H = 57584.57;  I = 0;        J = 57584.57;  K = 0;        L = 57584.57;  N = -67864.57000000001 }
    140 = @{ H = 35095;     I = 0;        J = 35095;     K = 0;        L = 35095;     N = -45455 }
    141 = @{ H = 0;         I = 0;        J = 0;         K = 0;        L = 0 }
}

foreach ($row in $bsmRows.Keys) {
    $cols = $bsmRows[$row]
    foreach ($col in $cols.Keys) {
        $wsBSM.Range("$col$row").Value = $cols[$col]
    }
}

# ---------------------------------------------------------------------------
# GSM sheet: existing H:N figures refreshed to new market-board averages.
# ---------------------------------------------------------------------------
$wsGSM = $wb.Worksheets.Item("GSM")

$gsmRows = @{
    70 = @{ H = 6059.8623;  I = 5875.9375; J = 6286.231;  K = 5875.9375; L = 6286.231;  M = -5605.9375; N = -6826.231 }
    73 = @{ H = 6059.8623;  I = 5875.9375; J = 6286.231;  K = 5875.9375; L = 6286.231;  M = -4939.9375; N = -8158.231 }
    80 = @{ H = 12431.044;  I = 6326;      J = 15687.066; K = 6326;      L = 15687.066; M = -5328;      N = -17683.066 }
    83 = @{ H = 12431.044;  I = 6326;      J = 15687.066; K = 31630;     L = 78435.33;  M = -26638;     N = -88419.33 }
}

foreach ($row in $gsmRows.Keys) {
    $cols = $gsmRows[$row]
    foreach ($col in $cols.Keys) {
        $wsGSM.Range("$col$row").Value = $cols[$col]
    }
}
